$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Plan1")
$ws2 = $wb.Worksheets.Item("Plan1 (2)")

# --- Plan1 (2): add new "explosion sequence" test cells to row 2 (P2:T2) ---
$ws2.Range("P2").Value = 0
$ws2.Range("Q2").Value = 43
$ws2.Range("R2").Formula = "=Q2+43"
$ws2.Range("S2:T2").Formula = "=R2+43"

# --- Update selections on both sheets ---
$ws1.Activate() | Out-Null
$ws1.Range("Q8").Select() | Out-Null

# Plan1 (2) becomes the active/selected sheet, with C3 selected
$ws2.Activate() | Out-Null
$ws2.Range("C3").Select() | Out-Null
